$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("D9 ->")

# Update "CURRENT AS OF" period header from D1.1 to D1.2
$ws.Range("T1").Value = "D1.2"

# Record D1.2 expenditures for the relevant ordnance rows (column E = D1.2)
$ws.Range("E10").Value = 4
$ws.Range("E11").Value = 2
$ws.Range("E12").Value = 4
$ws.Range("E15").Value = 5
$ws.Range("E18").Value = 5
$ws.Range("E19").Value = 12

# Update the active cell selection
$ws.Activate() | Out-Null
$ws.Range("H37").Select() | Out-Null
